$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 33: New York state hospitalization data for 16 April 2020.
# Copy the row above first so the new row inherits the same cell styles
# (date format in column A, general format elsewhere).
$ws.Range("A32:G32").Copy()
$ws.Range("A33:G33").PasteSpecial(-4122)

$ws.Range("A33").Value = 43937
$ws.Range("B33").Value = -419
$ws.Range("C33").Value = -32
$ws.Range("D33").Value = -73
$ws.Range("E33").ClearContents()
$ws.Range("F33").Value = 630
$ws.Range("G33").Value = 1974

# Update the selection to reflect the new active cell after the edit
$ws.Range("F34").Select()
